$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.851.60'
$ws.Range('E2').Value = '  -3.77%  '
$ws.Range('D3').Value = '1.667.41'
$ws.Range('E3').Value = '  -3.93%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.007'
$ws.Range('E4').Value = '  +0.78%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '309.30'
$ws.Range('E5').Value = '  -2.19%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.004'
$ws.Range('E6').Value = '  +1.07%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.3655'
$ws.Range('E7').Value = '  -4.18%  '
$ws.Range('B8').Value = 'Cardano'
$ws.Range('C8').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3332'
$ws.Range('E8').Value = '  -8.32%  '
$ws.Range('B9').Value = 'OKB'
$ws.Range('C9').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '47.11'
$ws.Range('E9').Value = '  -6.96%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.156'
$ws.Range('E10').Value = '  -5.56%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07225'
$ws.Range('E11').Value = '  -6.07%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.003'
$ws.Range('E12').Value = '  +0.78%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '6.116'
$ws.Range('E13').Value = '  -5.13%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '20.22'
$ws.Range('E14').Value = '  -7.01%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '6.736'
$ws.Range('E15').Value = '  -4.61%  '
$ws.Range('D16').Value = '1.663.97'
$ws.Range('E16').Value = '  -4.32%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.00001086'
$ws.Range('E17').Value = '  -5.89%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '1.003'
$ws.Range('E18').Value = '  +1.07%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06617'
$ws.Range('E19').Value = '  -2.96%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '81.10'
$ws.Range('E20').Value = '  -6.80%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '16.61'
$ws.Range('E21').Value = '  -5.36%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.094'
$ws.Range('E22').Value = '  -5.67%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '12.44'
$ws.Range('E23').Value = '  -2.58%  '
$ws.Range('D24').Value = '24.808.57'
$ws.Range('E24').Value = '  -3.71%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.427'
$ws.Range('E25').Value = '  -0.01%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.637'
$ws.Range('E26').Value = '  -10.02%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '149.04'
$ws.Range('E27').Value = '  -3.38%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '19.60'
$ws.Range('E28').Value = '  -5.01%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '128.60'
$ws.Range('E29').Value = '  -4.14%  '
$ws.Range('D30').Value = '1.850.55'
$ws.Range('E30').Value = '  -4.54%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.199'
$ws.Range('E31').Value = '  +0.26%  '
$ws.Range('B32').Value = 'HuobiToken'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.148'
$ws.Range('E32').Value = '  -5.15%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '6.335'
$ws.Range('E33').Value = '  -9.72%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.732'
$ws.Range('E34').Value = '  -3.68%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.08528'
$ws.Range('E35').Value = '  -1.73%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '13.18'
$ws.Range('E36').Value = '  -7.60%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '5.333'
$ws.Range('E37').Value = '  -5.39%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.06374'
$ws.Range('E38').Value = '  -5.18%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.02316'
$ws.Range('E39').Value = '  -5.83%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '8.628'
$ws.Range('E40').Value = '  -7.06%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.250'
$ws.Range('E41').Value = '  -3.66%  '
$ws.Range('B42').Value = 'Algorand'
$ws.Range('C42').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.2127'
$ws.Range('E42').Value = '  -3.94%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.6174'
$ws.Range('E43').Value = '  -5.64%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.003'
$ws.Range('E44').Value = '  +1.07%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '13.22'
$ws.Range('E45').Value = '  -5.09%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '3.784'
$ws.Range('E46').Value = '  -2.98%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.5877'
$ws.Range('E47').Value = '  -7.25%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.015'
$ws.Range('E48').Value = '  -7.28%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '124.73'
$ws.Range('E49').Value = '  -5.38%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.07100'
$ws.Range('E50').Value = '  -4.71%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '76.07'
$ws.Range('E51').Value = '  -4.01%  '
